$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.515984
$ws.Range("H2").Value = 7.547952
$ws.Range("I2").Value = 0.08781336966822693
$ws.Range("J2").Value = 0.09884082726736673
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.515984
$ws.Range("N2").Value = 7.547952
$ws.Range("O2").Value = 0.08781336966822693
$ws.Range("P2").Value = 0.09884082726736673
$ws.Range("Q2").Value = 6.330175488256
$ws.Range("R2").Value = 56.97157939430399
$ws.Range("S2").Value = 0.007711187892488678
$ws.Range("T2").Value = 0.009769509134897427

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.515984
$ws.Range("H3").Value = 7.547952
$ws.Range("I3").Value = 0.08781336966822693
$ws.Range("J3").Value = 0.09884082726736673
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.10016866666667
$ws.Range("N3").Value = 30.300506
$ws.Range("O3").Value = 0.3525180783492434
$ws.Range("P3").Value = 0.3967867150797739
$ws.Range("Q3").Value = 25.41186276263467
$ws.Range("R3").Value = 228.706764863712
$ws.Range("S3").Value = 0.0309558003288151
$ws.Range("T3").Value = 0.03921872716718579

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.515984
$ws.Range("H4").Value = 7.547952
$ws.Range("I4").Value = 0.08781336966822693
$ws.Range("J4").Value = 0.09884082726736673
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.445583666666667
$ws.Range("N4").Value = 19.336751
$ws.Range("O4").Value = 0.2249650320703493
$ws.Range("P4").Value = 0.2532157683969216
$ws.Range("Q4").Value = 16.21698537599467
$ws.Range("R4").Value = 145.952868383952
$ws.Range("S4").Value = 0.01975493752361811
$ws.Range("T4").Value = 0.02502805602549367

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.515984
$ws.Range("H5").Value = 7.547952
$ws.Range("I5").Value = 0.08781336966822693
$ws.Range("J5").Value = 0.09884082726736673
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 9.589755
$ws.Range("N5").Value = 19.17951
$ws.Range("O5").Value = 0.3347035199121805
$ws.Range("P5").Value = 0.2511566892559376
$ws.Range("Q5").Value = 24.12767014392
$ws.Range("R5").Value = 144.76602086352
$ws.Range("S5").Value = 0.02939144392330506
$ws.Range("T5").Value = 0.02482453493978984

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 10.10016866666667
$ws.Range("H6").Value = 30.300506
$ws.Range("I6").Value = 0.3525180783492434
$ws.Range("J6").Value = 0.3967867150797739
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.515984
$ws.Range("N6").Value = 7.547952
$ws.Range("O6").Value = 0.08781336966822693
$ws.Range("P6").Value = 0.09884082726736673
$ws.Range("Q6").Value = 25.41186276263467
$ws.Range("R6").Value = 228.706764863712
$ws.Range("S6").Value = 0.0309558003288151
$ws.Range("T6").Value = 0.03921872716718579

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 10.10016866666667
$ws.Range("H7").Value = 30.300506
$ws.Range("I7").Value = 0.3525180783492434
$ws.Range("J7").Value = 0.3967867150797739
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.10016866666667
$ws.Range("N7").Value = 30.300506
$ws.Range("O7").Value = 0.3525180783492434
$ws.Range("P7").Value = 0.3967867150797739
$ws.Range("Q7").Value = 102.0134070951151
$ws.Range("R7").Value = 918.1206638560359
$ws.Range("S7").Value = 0.1242689955630433
$ws.Range("T7").Value = 0.1574396972637976

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 10.10016866666667
$ws.Range("H8").Value = 30.300506
$ws.Range("I8").Value = 0.3525180783492434
$ws.Range("J8").Value = 0.3967867150797739
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 6.445583666666667
$ws.Range("N8").Value = 19.336751
$ws.Range("O8").Value = 0.2249650320703493
$ws.Range("P8").Value = 0.2532157683969216
$ws.Range("Q8").Value = 65.10148218844512
$ws.Range("R8").Value = 585.913339696006
$ws.Range("S8").Value = 0.07930424080121543
$ws.Range("T8").Value = 0.1004726529486153

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 10.10016866666667
$ws.Range("H9").Value = 30.300506
$ws.Range("I9").Value = 0.3525180783492434
$ws.Range("J9").Value = 0.3967867150797739
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 9.589755
$ws.Range("N9").Value = 19.17951
$ws.Range("O9").Value = 0.3347035199121805
$ws.Range("P9").Value = 0.2511566892559376
$ws.Range("Q9").Value = 96.85814297201
$ws.Range("R9").Value = 581.14885783206
$ws.Range("S9").Value = 0.1179890416561696
$ws.Range("T9").Value = 0.09965563770017503

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 6.445583666666667
$ws.Range("H10").Value = 19.336751
$ws.Range("I10").Value = 0.2249650320703493
$ws.Range("J10").Value = 0.2532157683969216
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.515984
$ws.Range("N10").Value = 7.547952
$ws.Range("O10").Value = 0.08781336966822693
$ws.Range("P10").Value = 0.09884082726736673
$ws.Range("Q10").Value = 16.21698537599467
$ws.Range("R10").Value = 145.952868383952
$ws.Range("S10").Value = 0.01975493752361811
$ws.Range("T10").Value = 0.02502805602549367

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 6.445583666666667
$ws.Range("H11").Value = 19.336751
$ws.Range("I11").Value = 0.2249650320703493
$ws.Range("J11").Value = 0.2532157683969216
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 10.10016866666667
$ws.Range("N11").Value = 30.300506
$ws.Range("O11").Value = 0.3525180783492434
$ws.Range("P11").Value = 0.3967867150797739
$ws.Range("Q11").Value = 65.10148218844512
$ws.Range("R11").Value = 585.913339696006
$ws.Range("S11").Value = 0.07930424080121543
$ws.Range("T11").Value = 0.1004726529486153

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 6.445583666666667
$ws.Range("H12").Value = 19.336751
$ws.Range("I12").Value = 0.2249650320703493
$ws.Range("J12").Value = 0.2532157683969216
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 6.445583666666667
$ws.Range("N12").Value = 19.336751
$ws.Range("O12").Value = 0.2249650320703493
$ws.Range("P12").Value = 0.2532157683969216
$ws.Range("Q12").Value = 41.54554880400011
$ws.Range("R12").Value = 373.909939236001
$ws.Range("S12").Value = 0.05060926565441327
$ws.Range("T12").Value = 0.06411822536484345

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 6.445583666666667
$ws.Range("H13").Value = 19.336751
$ws.Range("I13").Value = 0.2249650320703493
$ws.Range("J13").Value = 0.2532157683969216
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 9.589755
$ws.Range("N13").Value = 19.17951
$ws.Range("O13").Value = 0.3347035199121805
$ws.Range("P13").Value = 0.2511566892559376
$ws.Range("Q13").Value = 61.81156819533501
$ws.Range("R13").Value = 370.86940917201
$ws.Range("S13").Value = 0.07529658809110247
$ws.Range("T13").Value = 0.06359683405796912

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 9.589755
$ws.Range("H14").Value = 19.17951
$ws.Range("I14").Value = 0.3347035199121805
$ws.Range("J14").Value = 0.2511566892559376
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.515984
$ws.Range("N14").Value = 7.547952
$ws.Range("O14").Value = 0.08781336966822693
$ws.Range("P14").Value = 0.09884082726736673
$ws.Range("Q14").Value = 24.12767014392
$ws.Range("R14").Value = 144.76602086352
$ws.Range("S14").Value = 0.02939144392330506
$ws.Range("T14").Value = 0.02482453493978984

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 9.589755
$ws.Range("H15").Value = 19.17951
$ws.Range("I15").Value = 0.3347035199121805
$ws.Range("J15").Value = 0.2511566892559376
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 10.10016866666667
$ws.Range("N15").Value = 30.300506
$ws.Range("O15").Value = 0.3525180783492434
$ws.Range("P15").Value = 0.3967867150797739
$ws.Range("Q15").Value = 96.85814297201
$ws.Range("R15").Value = 581.14885783206
$ws.Range("S15").Value = 0.1179890416561696
$ws.Range("T15").Value = 0.09965563770017503

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 9.589755
$ws.Range("H16").Value = 19.17951
$ws.Range("I16").Value = 0.3347035199121805
$ws.Range("J16").Value = 0.2511566892559376
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 6.445583666666667
$ws.Range("N16").Value = 19.336751
$ws.Range("O16").Value = 0.2249650320703493
$ws.Range("P16").Value = 0.2532157683969216
$ws.Range("Q16").Value = 61.81156819533501
$ws.Range("R16").Value = 370.86940917201
$ws.Range("S16").Value = 0.07529658809110247
$ws.Range("T16").Value = 0.06359683405796912

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 9.589755
$ws.Range("H17").Value = 19.17951
$ws.Range("I17").Value = 0.3347035199121805
$ws.Range("J17").Value = 0.2511566892559376
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 9.589755
$ws.Range("N17").Value = 19.17951
$ws.Range("O17").Value = 0.2249650320703493
$ws.Range("P17").Value = 0.2511566892559376
$ws.Range("Q17").Value = 91.963400960025
$ws.Range("R17").Value = 367.8536038401
$ws.Range("S17").Value = 0.1120264462416034
$ws.Range("T17").Value = 0.06307968255800363

Write-Host "Updated rows 2-17"